# Apply cryptocurrency price/volume updates to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.142.27"
$ws.Range("E2").Value = "  +3.45%  "

# Row 3
$ws.Range("D3").Value = "1.599.97"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.49%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.16%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.485"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.99%  "

# Row 8
$ws.Range("E8").Value = "  +2.56%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0617"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.83%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.10%  "

# Row 12
$ws.Range("D12").Value = "1.823.20"
$ws.Range("E12").Value = "  +2.36%  "

# Row 13
$ws.Range("D13").Value = "1.607.28"
$ws.Range("E13").Value = "  +1.71%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.31%  "

# Row 16
$ws.Range("D16").Value = "26.109.25"
$ws.Range("E16").Value = "  +3.34%  "

# Row 17
$ws.Range("E17").Value = "  +2.09%  "

# Row 18
$ws.Range("D18").Value = "0.0₃0722"
$ws.Range("E18").Value = "  +1.37%  "

# Row 19
$ws.Range("E19").Value = "  -0.12%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "206.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +11.20%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.80%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.62%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.94%  "

# Row 24
$ws.Range("E24").Value = "  +11.78%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.57%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.125"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.08%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.68%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.44"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.16%  "

# Row 30
$ws.Range("E30").Value = "  +1.78%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0470"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.47%  "

# Row 32
$ws.Range("E32").Value = "  +2.91%  "

# Row 33
$ws.Range("E33").Value = "  +0.16%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.27%  "

# Row 35
$ws.Range("E35").Value = "  +2.39%  "

# Row 36
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0162"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.74%  "

# Row 37
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "1.108.09"
$ws.Range("E37").Value = "  +1.99%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.26%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.32"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.09%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.779"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.03%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.494"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.11%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.778"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.21%  "

# Row 43
$ws.Range("D43").Value = "1.740.14"
$ws.Range("E43").Value = "  +2.63%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "92.69"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.66%  "

# Row 45
$ws.Range("E45").Value = "  +0.77%  "

# Row 46
$ws.Range("D46").Value = "0.0₆0107"
$ws.Range("E46").Value = "  -3.41%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.04%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.47%  "

# Row 49
$ws.Range("E49").Value = "  +0.25%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.409"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.67%  "

# Row 51
$ws.Range("E51").Value = "  +0.04%  "
